$d = $word.ActiveDocument

# --- Update the date line at the top of the document ---
$d.Content.Find.Execute("2024-06-18 Tuesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-06-19 Wednesday", 2)

# --- Update every arithmetic expression in the practice table ---
# The table is 20 rows x 5 columns = 100 cells, filled row-major (left to
# right, top to bottom) with a single "a+b=" / "a-b=" expression per cell.
# This ordered list supplies the new text for each cell in that same
# row-major order.
$newValues = @(
    "35+26=", "90-4=", "47+0=", "32+31=", "85-11=", "3+15=", "41+47=", "25+31=", "41+43=", "88-86=",
    "75+4=", "9+7=", "60+37=", "81+15=", "81-80=", "1+94=", "44+50=", "23+53=", "68-61=", "92-50=",
    "73-13=", "11+64=", "92-44=", "8+60=", "76-66=", "51-22=", "37+12=", "83+1=", "44-26=", "56-33=",
    "81+3=", "72-64=", "94-71=", "9+84=", "22+36=", "62+34=", "78-24=", "87-11=", "64+19=", "71-16=",
    "76-75=", "77-39=", "47-47=", "98-38=", "37+2=", "47+37=", "47-24=", "5+2=", "52+39=", "5+77=",
    "92-88=", "7-3=", "59+19=", "11+79=", "17+43=", "24+42=", "68-52=", "79-66=", "29-25=", "46-4=",
    "8+38=", "82-13=", "70-36=", "5+85=", "26+32=", "36+9=", "80-13=", "5+53=", "34-9=", "85-18=",
    "90-81=", "76+17=", "30+43=", "32+32=", "19+33=", "2+29=", "95-84=", "91-53=", "88-39=", "69-58=",
    "96-74=", "19+80=", "40+13=", "88-54=", "6+92=", "99-24=", "8+43=", "64+25=", "31+22=", "76-59=",
    "77+10=", "13+16=", "20+36=", "34+3=", "13+49=", "17+65=", "9+3=", "24+26=", "22-22=", "60-12="
)

$tbl = $d.Tables.Item(1)
$rows = $tbl.Rows.Count
$cols = $tbl.Columns.Count

$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $tbl.Cell($r, $c)
        $cell.Range.Text = $newValues[$idx]
        $idx = $idx + 1
    }
}
